$d = $word.ActiveDocument

function Get-ParaIndexContainingText($searchText) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $f.Text = $searchText
    $f.Execute() | Out-Null
    if (-not $f.Found) {
        throw "Text not found: $searchText"
    }
    $startPos = $f.Parent.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $startPos -and $p.Range.End -gt $startPos) {
            return $i
        }
    }
    throw "Paragraph not found for text: $searchText"
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) "C;ip;nickname" paragraph: drop the underline char-formatting from the
#    paragraph mark, drop the _GoBack bookmark, and wrap "C;ip" in a
#    gramStart/gramEnd proofErr pair (in addition to the existing
#    spellStart/spellEnd that already wraps the whole run sequence).
# ---------------------------------------------------------------------------

$cipIdx = Get-ParaIndexContainingText("C;ip;nickname")
$cipRange = $d.Paragraphs.Item($cipIdx).Range

$cipXml = $pkgHeader + '<w:body><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:i/></w:rPr></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:i/></w:rPr><w:t>C;ip</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:rPr><w:i/></w:rPr><w:t>;nickname</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p></w:body>' + $pkgFooter

$cipRange.InsertXML($cipXml)

# ---------------------------------------------------------------------------
# 2) Insert two brand-new paragraphs ("LIN;..." and "LOUT;...") right after
#    the "P;(lettera1)-...;punteggioParola" paragraph, before "Disconnessione:".
#    The second one carries the _GoBack bookmark that used to live on the
#    "C;ip;nickname" paragraph.
# ---------------------------------------------------------------------------

$discIdx = Get-ParaIndexContainingText("Disconnessione:")
$insertPoint = $d.Paragraphs.Item($discIdx).Range
$insertPoint.Collapse(1)

$linLoutXml = $pkgHeader + '<w:body>' + `
    '<w:p><w:pPr><w:ind w:left="360"/></w:pPr>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>L</w:t></w:r>' + `
        '<w:r><w:t>IN</w:t></w:r>' + `
        '<w:r><w:t>;lettera</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t>1;lettera2;&#8230;;</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
        '<w:r><w:t>sempre</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> 8 lettere) </w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:ind w:left="360"/></w:pPr>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>L</w:t></w:r>' + `
        '<w:r><w:t>OUT</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
        '<w:bookmarkEnd w:id="0"/>' + `
        '<w:r><w:t>;lettera</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t>1;lettera2;..;</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>max</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> 8 lettere)</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body>' + $pkgFooter

$insertPoint.InsertXML($linLoutXml)
